$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "241027"
$ws.Range("F6").Value = "241127"

$ws.Range("F7").Select()
